$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 112
$ws1.Range("F3").Value = 7439
$ws1.Range("F4").Value = 278
$ws1.Range("F6").Value = 3999
$ws1.Range("F9").Value = 273
$ws1.Range("F11").Value = 122

# Sheet "演出" (sheet2): update "想去人数" (column F) counts
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 4

# Sheet "全部类型" (sheet4): update "想去人数" (column F) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 112
$ws4.Range("F4").Value = 7439
$ws4.Range("F6").Value = 278
$ws4.Range("F8").Value = 3999
$ws4.Range("F11").Value = 273
$ws4.Range("F13").Value = 4
$ws4.Range("F14").Value = 122

$wb.Save()
